$d = $word.ActiveDocument

$replacements = @(
    @("728÷5=", "471÷5="),
    @("867÷8=", "632÷4="),
    @("243÷4=", "443÷9="),
    @("331÷2=", "571÷8="),
    @("731÷8=", "246÷6="),
    @("737÷6=", "722÷2="),
    @("186÷6=", "549÷5="),
    @("866÷4=", "226÷9="),
    @("648÷6=", "189÷8="),
    @("542÷7=", "251÷9="),
    @("381÷9=", "970÷8="),
    @("551÷9=", "441÷7="),
    @("641÷7=", "363÷8="),
    @("939÷8=", "646÷2="),
    @("319÷9=", "950÷6="),
    @("927÷9=", "514÷8="),
    @("834÷8=", "469÷4="),
    @("749÷2=", "846÷7="),
    @("570÷9=", "806÷3="),
    @("845÷7=", "668÷2="),
    @("153÷7=", "187÷4="),
    @("651÷8=", "455÷3="),
    @("592÷2=", "298÷6="),
    @("124÷7=", "930÷5="),
    @("709÷8=", "668÷8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
